$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.136.86"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "1.871.15"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5062"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3748"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07153"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8895"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "1.863.47"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07559"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.324"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008498"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "27.182.98"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.078"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "2.112.35"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.486"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.842"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.096"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.762"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.688"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09007"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05132"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.085"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7429"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("E36").Value = "  -5.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02039"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.534"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.043"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.075"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5388"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.597"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.447"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1479"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4646"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.571"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
